$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("V1").Value = 'CIN2+ SE immediate'
$ws.Range("W1").Value = 'CIN2+ LL95 immediate'
$ws.Range("X1").Value = 'CIN2+ UL95 immediate'
$ws.Range("Z1").Value = 'CIN2+ SE 1-year'
$ws.Range("AA1").Value = 'CIN2+ LL95 1-year'
$ws.Range("AB1").Value = 'CIN2+ UL95 1-year'
$ws.Range("AD1").Value = 'CIN2+ SE 2-year'
$ws.Range("AE1").Value = 'CIN2+ LL95 2-year'
$ws.Range("AF1").Value = 'CIN2+ UL95 2-year'
$ws.Range("AH1").Value = 'CIN2+ SE 3-year'
$ws.Range("AI1").Value = 'CIN2+ LL95 3-year'
$ws.Range("AJ1").Value = 'CIN2+ UL95 3-year'
$ws.Range("AL1").Value = 'SCIN2+ E 4-year'
$ws.Range("AM1").Value = 'CIN2+ LL95 4-year'
$ws.Range("AN1").Value = 'CIN2+ UL95 4-year'
$ws.Range("AP1").Value = 'CIN2+ SE 5-year'
$ws.Range("AQ1").Value = 'CIN2+ LL95 5-year'
$ws.Range("AR1").Value = 'CIN2+ UL95 5-year'
$ws.Range("AT1").Value = 'CIN3+ SE immediate'
$ws.Range("AU1").Value = 'CIN3+ LL95 immediate'
$ws.Range("AV1").Value = 'CIN3+ UL95 immediate'
$ws.Range("AX1").Value = 'CIN3+ SE 1-year'
$ws.Range("AY1").Value = 'CIN3+ LL95 1-year'
$ws.Range("AZ1").Value = 'CIN3+ UL95 1-year'
$ws.Range("BB1").Value = 'CIN3+ SE 2-year'
$ws.Range("BC1").Value = 'CIN3+ LL95 2-year'
$ws.Range("BD1").Value = 'CIN3+ UL95 2-year'
$ws.Range("BF1").Value = 'CIN3+ SE 3-year'
$ws.Range("BG1").Value = 'CIN3+ L95 3-year'
$ws.Range("BH1").Value = 'CIN3+ UL95 3-year'
$ws.Range("BJ1").Value = 'CIN3+ SE 4-year'
$ws.Range("BK1").Value = 'CIN3+ LL95 4-year'
$ws.Range("BL1").Value = 'CIN3+ UL95 4-year'
$ws.Range("BN1").Value = 'CIN3+ SE 5-year'
$ws.Range("BO1").Value = 'CIN3+ LL95 5-year'
$ws.Range("BP1").Value = 'CIN3+ UL95 5-year'
$ws.Range("BR1").Value = 'CANCER SE immediate'
$ws.Range("BS1").Value = 'CANCER LL95 immediate'
$ws.Range("BT1").Value = 'CANCER UL95 immediate'
$ws.Range("BV1").Value = 'CANCER SE 1-year'
$ws.Range("BW1").Value = 'CANCER LL95 1-year'
$ws.Range("BX1").Value = 'CANCER UL95 1-year'
$ws.Range("BZ1").Value = 'CANCER SE 2-year'
$ws.Range("CA1").Value = 'CANCER LL95 2-year'
$ws.Range("CB1").Value = 'CANCER UL95 2-year'
$ws.Range("CD1").Value = 'CANCER SE 3-year'
$ws.Range("CE1").Value = 'CANCER LL95 3-year'
$ws.Range("CF1").Value = 'CANCER UL95 3-year'
$ws.Range("CH1").Value = 'CANCER SE 4-year'
$ws.Range("CI1").Value = 'CANCER LL95 4-year'
$ws.Range("CJ1").Value = 'CANCER UL95 4-year'
$ws.Range("CL1").Value = 'CANCER SE 5-year'
$ws.Range("CM1").Value = 'CANCER LL95 5-year'
$ws.Range("CN1").Value = 'CANCER UL95 5-year'
$ws.Range("CS1").Value = 'Unweighted %'
$ws.Range("CT1").Value = 'Unweighted Informative N'
$ws.Range("CV1").Value = 'Unweighted CIN2+ Prevalence Cases'
$ws.Range("CW1").Value = 'Unweighted CIN2+ Incidence Cases'
$ws.Range("CX1").Value = 'Unweighted CIN2+ Unknown Cases'
$ws.Range("CY1").Value = 'Unweighted Number of CIN3+ Cases'
$ws.Range("CZ1").Value = 'Unweighted CIN3+ Prevalence Cases'
$ws.Range("DA1").Value = 'Unweighted CIN3+ Incidence Cases'
$ws.Range("DB1").Value = 'Unweighted CIN3+ Unknown Cases'
$ws.Range("DD1").Value = 'Unweighted Cancer Prevalence Cases'
$ws.Range("DE1").Value = 'Unweighted Cancer Incidence Cases'
$ws.Range("DF1").Value = 'Unweighted Cancer Unknown Cases'
